$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E (boson), shifting every following
# column one place to the right (E->F, F->G, ... X->Y).
$ws.Columns("E").Insert()

# Populate the new "pt_max" column: header in row 1 and a constant value
# of 50 for every data row (2-23).
$ws.Range("E1").Value = "pt_max"
$ws.Range("E2:E23").Value = 50

# The insert operation can leave the very last shared-formula group of a
# row (now column U) serialized as individual, non-shared formulas.
# Re-apply the formula across the whole range so it is written back out
# as a single shared formula group, matching the expected output.
$ws.Range("U2:U23").Formula = "=N2/100*G2"

# Restore the selection described by the change (E2:E23, active cell E2).
$ws.Range("E2:E23").Select()
